# Insert a new row at position 79, shifting existing rows 79-154 down to 80-155.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(79).Insert()

# Populate the newly inserted row 79 with the new record (same as the row
# that used to be there, except for the Date and Volumen columns).
$ws.Range("A79").Value = 10
$ws.Range("B79").Value = "Vega Modelo de Temuco"
$ws.Range("C79").Value = "La Araucanía"
$ws.Range("D79").Value = 44512
$ws.Range("E79").Value = 9
$ws.Range("F79").Value = 100112005
$ws.Range("G79").Value = "Puerro"
$ws.Range("H79").Value = "Azul de Maquehue"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 40
$ws.Range("K79").Value = 8000
$ws.Range("L79").Value = 8000
$ws.Range("M79").Value = 8000
$ws.Range("N79").Value = "$/docena de paquetes"
$ws.Range("O79").Value = "Provincia de Cautín"
$ws.Range("P79").Value = 667
$ws.Range("Q79").Value = 12
$ws.Range("R79").Value = "Hortaliza"
